$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score")

# Minute2 (col F) / Second2 (col G) are now populated for every athlete row (2-21)
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 6).Value = 44
    $ws.Cells.Item($r, 7).Value = 0
}

# Rep2a..Rep2f (cols H-M) plus the Rep2 total (col N, =SUM(H:M)) are filled in
# for the first ten athletes (rows 2-11); the remaining rows only got the
# Minute2/Second2 values set above.
$rep2Rows = 2,3,4,5,6,7,8,9,10,11
$rep2Data = @(
    (54, 52, 166, 134, 4.5,  171),
    (46, 49, 130, 131, 14,   128),
    (50, 56, 163, 143, 19,   168),
    (62, 57, 166, 155, 17.5, 202),
    (68, 57, 160, 158, 20.5, 180),
    (63, 53, 155, 159, 12,   162),
    (50, 32, 133, 125, 8,    162),
    (52, 49, 138, 135, 17,   149),
    (52, 57, 158, 139, 21.5, 174),
    (54, 49, 155, 144, 19,   160)
)

for ($i = 0; $i -lt $rep2Rows.Length; $i++) {
    $r = $rep2Rows[$i]
    $vals = $rep2Data[$i]
    $ws.Cells.Item($r, 8).Value  = $vals[0]
    $ws.Cells.Item($r, 9).Value  = $vals[1]
    $ws.Cells.Item($r, 10).Value = $vals[2]
    $ws.Cells.Item($r, 11).Value = $vals[3]
    $ws.Cells.Item($r, 12).Value = $vals[4]
    $ws.Cells.Item($r, 13).Value = $vals[5]
    $ws.Range("N$r").Formula = "=SUM(H" + $r + ":M" + $r + ")"
}

# Match the author's final cursor position
$ws.Range("P16").Select()
